$p = $ppt.ActivePresentation
try {
  $r = $p.Slides.InsertFromFile("/tmp/work/donor.pptx", $p.Slides.Count)
  Write-Output ("InsertFromFile result: " + $r)
} catch {
  Write-Output ("EXC: " + $_.Exception.Message)
}
Write-Output ("Slides.Count=" + $p.Slides.Count)
Write-Output ("IronPkgSize=" + $p.IronPkgSize)
